# Regenerate merged AHB files
# - Rename the "_old"/"_new" header-label suffixes to "_FV2404"/"_FV2410"
# - Wrap the data range in a native Excel Table (ListObject)
# - Freeze the header row (pane split under row 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$lastCol = $ws.UsedRange.Columns.Count

# --- Rename header cells (row 1) -----------------------------------------
for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Value
    if ($val -ne $null) {
        $newVal = $val.ToString()
        $newVal = $newVal.Replace("_old", "_FV2404")
        $newVal = $newVal.Replace("_new", "_FV2410")
        if ($newVal -ne $val) {
            $cell.Value = $newVal
        }
    }
}

# --- Wrap data range into a Table (ListObject) ---------------------------
$dataRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $dataRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

# --- Freeze the header row (split below row 1) ----------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
